$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# URL
$meta.Range("B2").Value = "https://department-of-veterans-affairs.github.io/mhv-fhir-phr-mapping/ValueSet/DocumentReferenceTypeVS"

# Version
$meta.Range("B3").Value = "0.2.0"

# Date
$meta.Range("B8").Value = "2023-08-22T16:36:15-05:00"

# Publisher
$meta.Range("B9").Value = "VA Digital Services"

# --- Include ValueSets sheet updates ---
$incVS = $wb.Worksheets.Item("Include ValueSets")

# ValueSet URL
$incVS.Range("A2").Value = "https://department-of-veterans-affairs.github.io/mhv-fhir-phr-mapping/ValueSet/NoteTypeVS"
